$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray duplicated header row that was accidentally left in the
# data (row 57 contains "Numero de departement" / "Nom du departement" / "Nom de la region").
$ws.Rows.Item(57).Delete()

# Re-sort the data (now rows 1:97 incl. header) by column A (numero_departement)
# ascending, instead of the previous sort by column C (region_departement).
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A2:A97"))
$sortObj.SetRange($ws.Range("A1:D97"))
$sortObj.Header = 1
$sortObj.Apply()

# Update the selection to match the latest saved view: the whole first row selected.
$ws.Rows.Item(1).Select()
